$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing shared string text ("FirstTweet" -> "First Tweet")
$ws.Range("B2").Value = "It's my First Tweet!"

# Add new rows 3 and 4 with additional tweet data
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "It's my Second Tweet!"
$ws.Range("C3").Formula = "=TRUE"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "It's my Third Tweet!"
$ws.Range("C4").Formula = "=TRUE"

# Widen column B to fit the new, longer text
$ws.Columns.Item(2).ColumnWidth = 19.7

# Update the active selection to C4, matching the new last-entered cell
$ws.Range("C4").Select() | Out-Null
